# Applies the "ctdc gender male script with db in qa" edit:
#  - updates the Neo4j query text stored in cell B2 of the "startup" sheet
#    (adds an `age`/`weight` WITH-alias, wraps Age & Weight in an integer-coalescing
#    CASE expression, and appends a new `Cohort` output column)
#  - refreshes the (wrap-text driven) row heights for rows 2-4 that change as a
#    result of the longer query text

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
  MATCH (f:file)-[*]->(c)
    WHERE f.file_format IN ["tif"]  
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS Cohort
'@

$ws.Range("B2").Value = $newQuery

# The query cell wraps text (style applied to column B/C), so growing the text
# changes the auto-computed row height for this row (and, due to shared wrap
# width, its neighbours). Match the resulting heights.
$ws.Rows.Item(2).RowHeight = 290
$ws.Rows.Item(3).RowHeight = 290
$ws.Rows.Item(4).RowHeight = 261
